$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> M1
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Col3a1"
$ws.Cells.Item(2,3).Value = "Mag"
$ws.Cells.Item(2,4).Value = "M1"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 20.5550805
$ws.Cells.Item(2,8).Value = 41.11016100000001
$ws.Cells.Item(2,9).Value = 0.01692986717097462
$ws.Cells.Item(2,10).Value = 0.01176776206024777
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.27644
$ws.Cells.Item(2,14).Value = 0.8293199999999999
$ws.Cells.Item(2,15).Value = 0.1940440920813295
$ws.Cells.Item(2,16).Value = 0.2208748168298663
$ws.Cells.Item(2,17).Value = 5.68224645342
$ws.Cells.Item(2,18).Value = 34.09347872052
$ws.Cells.Item(2,19).Value = 0.003285140704249277
$ws.Cells.Item(2,20).Value = 0.002599202289554676

# Row 3: ECs -> M2
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Col3a1"
$ws.Cells.Item(3,3).Value = "Mag"
$ws.Cells.Item(3,4).Value = "M2"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 20.5550805
$ws.Cells.Item(3,8).Value = 41.11016100000001
$ws.Cells.Item(3,9).Value = 0.01692986717097462
$ws.Cells.Item(3,10).Value = 0.01176776206024777
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.6290166666666667
$ws.Cells.Item(3,14).Value = 1.88705
$ws.Cells.Item(3,15).Value = 0.4415315004607062
$ws.Cells.Item(3,16).Value = 0.502582625643659
$ws.Cells.Item(3,17).Value = 12.929488219175
$ws.Cells.Item(3,18).Value = 77.57692931505001
$ws.Cells.Item(3,19).Value = 0.007475069654600877
$ws.Cells.Item(3,20).Value = 0.005914272754189158

# Row 4: ECs -> sCs
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Col3a1"
$ws.Cells.Item(4,3).Value = "Mag"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 20.5550805
$ws.Cells.Item(4,8).Value = 41.11016100000001
$ws.Cells.Item(4,9).Value = 0.01692986717097462
$ws.Cells.Item(4,10).Value = 0.01176776206024777
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.5191680000000001
$ws.Cells.Item(4,14).Value = 1.038336
$ws.Cells.Item(4,15).Value = 0.3644244074579644
$ws.Cells.Item(4,16).Value = 0.2765425575264748
$ws.Cells.Item(4,17).Value = 10.671540033024
$ws.Cells.Item(4,18).Value = 42.68616013209601
$ws.Cells.Item(4,19).Value = 0.006169656812124473
$ws.Cells.Item(4,20).Value = 0.003254287016503937

# Row 5: FAPs -> M1
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Col3a1"
$ws.Cells.Item(5,3).Value = "Mag"
$ws.Cells.Item(5,4).Value = "M1"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1042.616902666667
$ws.Cells.Item(5,8).Value = 3127.850708
$ws.Cells.Item(5,9).Value = 0.8587349328240113
$ws.Cells.Item(5,10).Value = 0.8953456273674414
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.27644
$ws.Cells.Item(5,14).Value = 0.8293199999999999
$ws.Cells.Item(5,15).Value = 0.1940440920813295
$ws.Cells.Item(5,16).Value = 0.2208748168298663
$ws.Cells.Item(5,17).Value = 288.2210165731733
$ws.Cells.Item(5,18).Value = 2593.98914915856
$ws.Cells.Item(5,19).Value = 0.1666324403783567
$ws.Cells.Item(5,20).Value = 0.1977593014442053

# Row 6: FAPs -> M2
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Col3a1"
$ws.Cells.Item(6,3).Value = "Mag"
$ws.Cells.Item(6,4).Value = "M2"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1042.616902666667
$ws.Cells.Item(6,8).Value = 3127.850708
$ws.Cells.Item(6,9).Value = 0.8587349328240113
$ws.Cells.Item(6,10).Value = 0.8953456273674414
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.6290166666666667
$ws.Cells.Item(6,14).Value = 1.88705
$ws.Cells.Item(6,15).Value = 0.4415315004607062
$ws.Cells.Item(6,16).Value = 0.502582625643659
$ws.Cells.Item(6,17).Value = 655.8234087257111
$ws.Cells.Item(6,18).Value = 5902.410678531401
$ws.Cells.Item(6,19).Value = 0.3791585233878094
$ws.Cells.Item(6,20).Value = 0.4499851562608977

# Row 7: FAPs -> sCs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Col3a1"
$ws.Cells.Item(7,3).Value = "Mag"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1042.616902666667
$ws.Cells.Item(7,8).Value = 3127.850708
$ws.Cells.Item(7,9).Value = 0.8587349328240113
$ws.Cells.Item(7,10).Value = 0.8953456273674414
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.5191680000000001
$ws.Cells.Item(7,14).Value = 1.038336
$ws.Cells.Item(7,15).Value = 0.3644244074579644
$ws.Cells.Item(7,16).Value = 0.2765425575264748
$ws.Cells.Item(7,17).Value = 541.293332123648
$ws.Cells.Item(7,18).Value = 3247.759992741888
$ws.Cells.Item(7,19).Value = 0.3129439690578452
$ws.Cells.Item(7,20).Value = 0.2476011696623384

# Row 8: M1 -> M1
$ws.Cells.Item(8,1).Value = "M1"
$ws.Cells.Item(8,2).Value = "Col3a1"
$ws.Cells.Item(8,3).Value = "Mag"
$ws.Cells.Item(8,4).Value = "M1"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.8184900000000001
$ws.Cells.Item(8,8).Value = 2.45547
$ws.Cells.Item(8,9).Value = 0.0006741363518751979
$ws.Cells.Item(8,10).Value = 0.0007028770017727877
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.27644
$ws.Cells.Item(8,14).Value = 0.8293199999999999
$ws.Cells.Item(8,15).Value = 0.1940440920813295
$ws.Cells.Item(8,16).Value = 0.2208748168298663
$ws.Cells.Item(8,17).Value = 0.2262633756
$ws.Cells.Item(8,18).Value = 2.0363703804
$ws.Cells.Item(8,19).Value = 0.0001308121763386424
$ws.Cells.Item(8,20).Value = 0.0001552478290204901

# Row 9: M1 -> M2
$ws.Cells.Item(9,1).Value = "M1"
$ws.Cells.Item(9,2).Value = "Col3a1"
$ws.Cells.Item(9,3).Value = "Mag"
$ws.Cells.Item(9,4).Value = "M2"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.8184900000000001
$ws.Cells.Item(9,8).Value = 2.45547
$ws.Cells.Item(9,9).Value = 0.0006741363518751979
$ws.Cells.Item(9,10).Value = 0.0007028770017727877
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.6290166666666667
$ws.Cells.Item(9,14).Value = 1.88705
$ws.Cells.Item(9,15).Value = 0.4415315004607062
$ws.Cells.Item(9,16).Value = 0.502582625643659
$ws.Cells.Item(9,17).Value = 0.5148438515
$ws.Cells.Item(9,18).Value = 4.6335946635
$ws.Cells.Item(9,19).Value = 0.0002976524349585627
$ws.Cells.Item(9,20).Value = 0.0003532537690555104

# Row 10: M1 -> sCs
$ws.Cells.Item(10,1).Value = "M1"
$ws.Cells.Item(10,2).Value = "Col3a1"
$ws.Cells.Item(10,3).Value = "Mag"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.8184900000000001
$ws.Cells.Item(10,8).Value = 2.45547
$ws.Cells.Item(10,9).Value = 0.0006741363518751979
$ws.Cells.Item(10,10).Value = 0.0007028770017727877
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.5191680000000001
$ws.Cells.Item(10,14).Value = 1.038336
$ws.Cells.Item(10,15).Value = 0.3644244074579644
$ws.Cells.Item(10,16).Value = 0.2765425575264748
$ws.Cells.Item(10,17).Value = 0.4249338163200001
$ws.Cells.Item(10,18).Value = 2.54960289792
$ws.Cells.Item(10,19).Value = 0.0002456717405779928
$ws.Cells.Item(10,20).Value = 0.0001943754036967873

# Row 11: M2 -> M1
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Col3a1"
$ws.Cells.Item(11,3).Value = "Mag"
$ws.Cells.Item(11,4).Value = "M1"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 1.314656666666667
$ws.Cells.Item(11,8).Value = 3.94397
$ws.Cells.Item(11,9).Value = 0.001082796184724401
$ws.Cells.Item(11,10).Value = 0.001128959347368048
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 0.6666666666666666
$ws.Cells.Item(11,13).Value = 0.27644
$ws.Cells.Item(11,14).Value = 0.8293199999999999
$ws.Cells.Item(11,15).Value = 0.1940440920813295
$ws.Cells.Item(11,16).Value = 0.2208748168298663
$ws.Cells.Item(11,17).Value = 0.3634236889333333
$ws.Cells.Item(11,18).Value = 3.2708132004
$ws.Cells.Item(11,19).Value = 0.0002101102025739739
$ws.Cells.Item(11,20).Value = 0.0002493586890582831

# Row 12: M2 -> M2
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Col3a1"
$ws.Cells.Item(12,3).Value = "Mag"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 1.314656666666667
$ws.Cells.Item(12,8).Value = 3.94397
$ws.Cells.Item(12,9).Value = 0.001082796184724401
$ws.Cells.Item(12,10).Value = 0.001128959347368048
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.6290166666666667
$ws.Cells.Item(12,14).Value = 1.88705
$ws.Cells.Item(12,15).Value = 0.4415315004607062
$ws.Cells.Item(12,16).Value = 0.502582625643659
$ws.Cells.Item(12,17).Value = 0.8269409542777779
$ws.Cells.Item(12,18).Value = 7.442468588500001
$ws.Cells.Item(12,19).Value = 0.0004780886241344927
$ws.Cells.Item(12,20).Value = 0.0005673953530451854

# Row 13: M2 -> sCs
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Col3a1"
$ws.Cells.Item(13,3).Value = "Mag"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 1.314656666666667
$ws.Cells.Item(13,8).Value = 3.94397
$ws.Cells.Item(13,9).Value = 0.001082796184724401
$ws.Cells.Item(13,10).Value = 0.001128959347368048
$ws.Cells.Item(13,11).Value = 2
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.5191680000000001
$ws.Cells.Item(13,14).Value = 1.038336
$ws.Cells.Item(13,15).Value = 0.3644244074579644
$ws.Cells.Item(13,16).Value = 0.2765425575264748
$ws.Cells.Item(13,17).Value = 0.6825276723200002
$ws.Cells.Item(13,18).Value = 4.095166033920001
$ws.Cells.Item(13,19).Value = 0.0003945973580159344
$ws.Cells.Item(13,20).Value = 0.00031220530526458

# Row 14: Neutro -> M1
$ws.Cells.Item(14,1).Value = "Neutro"
$ws.Cells.Item(14,2).Value = "Col3a1"
$ws.Cells.Item(14,3).Value = "Mag"
$ws.Cells.Item(14,4).Value = "M1"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 20.443657
$ws.Cells.Item(14,8).Value = 61.330971
$ws.Cells.Item(14,9).Value = 0.01683809496630118
$ws.Cells.Item(14,10).Value = 0.01755595833477656
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 0.6666666666666666
$ws.Cells.Item(14,13).Value = 0.27644
$ws.Cells.Item(14,14).Value = 0.8293199999999999
$ws.Cells.Item(14,15).Value = 0.1940440920813295
$ws.Cells.Item(14,16).Value = 0.2208748168298663
$ws.Cells.Item(14,17).Value = 5.651444541079999
$ws.Cells.Item(14,18).Value = 50.86300086972
$ws.Cells.Item(14,19).Value = 0.003267332850115116
$ws.Cells.Item(14,20).Value = 0.003877669081466536

# Row 15: Neutro -> M2
$ws.Cells.Item(15,1).Value = "Neutro"
$ws.Cells.Item(15,2).Value = "Col3a1"
$ws.Cells.Item(15,3).Value = "Mag"
$ws.Cells.Item(15,4).Value = "M2"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 20.443657
$ws.Cells.Item(15,8).Value = 61.330971
$ws.Cells.Item(15,9).Value = 0.01683809496630118
$ws.Cells.Item(15,10).Value = 0.01755595833477656
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 0.6290166666666667
$ws.Cells.Item(15,14).Value = 1.88705
$ws.Cells.Item(15,15).Value = 0.4415315004607062
$ws.Cells.Item(15,16).Value = 0.502582625643659
$ws.Cells.Item(15,17).Value = 12.85940098061667
$ws.Cells.Item(15,18).Value = 115.73460882555
$ws.Cells.Item(15,19).Value = 0.007434549335370823
$ws.Cells.Item(15,20).Value = 0.008823319635582682

# Row 16: Neutro -> sCs
$ws.Cells.Item(16,1).Value = "Neutro"
$ws.Cells.Item(16,2).Value = "Col3a1"
$ws.Cells.Item(16,3).Value = "Mag"
$ws.Cells.Item(16,4).Value = "sCs"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 20.443657
$ws.Cells.Item(16,8).Value = 61.330971
$ws.Cells.Item(16,9).Value = 0.01683809496630118
$ws.Cells.Item(16,10).Value = 0.01755595833477656
$ws.Cells.Item(16,11).Value = 2
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.5191680000000001
$ws.Cells.Item(16,14).Value = 1.038336
$ws.Cells.Item(16,15).Value = 0.3644244074579644
$ws.Cells.Item(16,16).Value = 0.2765425575264748
$ws.Cells.Item(16,17).Value = 10.613692517376
$ws.Cells.Item(16,18).Value = 63.682155104256
$ws.Cells.Item(16,19).Value = 0.006136212780815241
$ws.Cells.Item(16,20).Value = 0.004854969617727341

# Row 17: sCs -> M1
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Col3a1"
$ws.Cells.Item(17,3).Value = "Mag"
$ws.Cells.Item(17,4).Value = "M1"
$ws.Cells.Item(17,5).Value = 2
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 128.382446
$ws.Cells.Item(17,8).Value = 256.764892
$ws.Cells.Item(17,9).Value = 0.1057401725021131
$ws.Cells.Item(17,10).Value = 0.07349881588839352
$ws.Cells.Item(17,11).Value = 2
$ws.Cells.Item(17,12).Value = 0.6666666666666666
$ws.Cells.Item(17,13).Value = 0.27644
$ws.Cells.Item(17,14).Value = 0.8293199999999999
$ws.Cells.Item(17,15).Value = 0.1940440920813295
$ws.Cells.Item(17,16).Value = 0.2208748168298663
$ws.Cells.Item(17,17).Value = 35.49004337224
$ws.Cells.Item(17,18).Value = 212.94026023344
$ws.Cells.Item(17,19).Value = 0.02051825576969571
$ws.Cells.Item(17,20).Value = 0.01623403749656099

# Row 18: sCs -> M2
$ws.Cells.Item(18,1).Value = "sCs"
$ws.Cells.Item(18,2).Value = "Col3a1"
$ws.Cells.Item(18,3).Value = "Mag"
$ws.Cells.Item(18,4).Value = "M2"
$ws.Cells.Item(18,5).Value = 2
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = 128.382446
$ws.Cells.Item(18,8).Value = 256.764892
$ws.Cells.Item(18,9).Value = 0.1057401725021131
$ws.Cells.Item(18,10).Value = 0.07349881588839352
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = 0.6290166666666667
$ws.Cells.Item(18,14).Value = 1.88705
$ws.Cells.Item(18,15).Value = 0.4415315004607062
$ws.Cells.Item(18,16).Value = 0.502582625643659
$ws.Cells.Item(18,17).Value = 80.75469824143333
$ws.Cells.Item(18,18).Value = 484.5281894486
$ws.Cells.Item(18,19).Value = 0.04668761702383192
$ws.Cells.Item(18,20).Value = 0.03693922787088869

# Row 19: sCs -> sCs
$ws.Cells.Item(19,1).Value = "sCs"
$ws.Cells.Item(19,2).Value = "Col3a1"
$ws.Cells.Item(19,3).Value = "Mag"
$ws.Cells.Item(19,4).Value = "sCs"
$ws.Cells.Item(19,5).Value = 2
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = 128.382446
$ws.Cells.Item(19,8).Value = 256.764892
$ws.Cells.Item(19,9).Value = 0.1057401725021131
$ws.Cells.Item(19,10).Value = 0.07349881588839352
$ws.Cells.Item(19,11).Value = 2
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = 0.5191680000000001
$ws.Cells.Item(19,14).Value = 1.038336
$ws.Cells.Item(19,15).Value = 0.3644244074579644
$ws.Cells.Item(19,16).Value = 0.2765425575264748
$ws.Cells.Item(19,17).Value = 66.652057724928
$ws.Cells.Item(19,18).Value = 266.608230899712
$ws.Cells.Item(19,19).Value = 0.03853429970858552
$ws.Cells.Item(19,20).Value = 0.02032555052094385
